# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets, which contain duplicated event data.
#
# Row -> new F value
#   5  -> 2811
#   11 -> 71
#   13 -> 1255
#   15 -> 389
#   17 -> 52
#   18 -> 46
#   22 -> 2779

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    5  = 2811
    11 = 71
    13 = 1255
    15 = 389
    17 = 52
    18 = 46
    22 = 2779
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
